# Add a new "2022-Q3" sheet (with fund holding detail) positioned right
# after the "总计" (total) sheet and before the existing "2022-Q2" sheet,
# and record its summary numbers on the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet before the current "2022-Q2"
#    sheet (index 2), so the final tab order is:
#    总计, 2022-Q3, 2022-Q2, 2022-Q1
# ---------------------------------------------------------------------
$wsQ2Before = $wb.Worksheets.Item(2)
$wsNew = $wb.Worksheets.Add($wsQ2Before)
$wsNew.Name = "2022-Q3"

# Re-resolve the "2022-Q2" sheet reference now that it has shifted to
# index 3 (references in this host track slot position, not identity).
$wsQ2 = $wb.Worksheets.Item(3)

# Seed the new sheet's header row + first data row with the same
# layout/formatting used by the "2022-Q2" sheet (column headers plus the
# thin-border / bold / centered style used on that sheet), then stamp a
# second data row using the same row-2 formatting.
$wsQ2.Range("B1:H1").Copy($wsNew.Range("B1"))
$wsQ2.Range("A2:H2").Copy($wsNew.Range("A2"))
$wsQ2.Range("A2:H2").Copy($wsNew.Range("A3"))

# ---------------------------------------------------------------------
# 2. Fill in the "2022-Q3" fund holding data.
# ---------------------------------------------------------------------
$wsNew.Range("A2").Value = 0
$wsNew.Range("B2").Value = "'011729"
$wsNew.Range("C2").Value = "工银聚享混合A"
$wsNew.Range("D2").Value = "'1.36"
$wsNew.Range("E2").Value = "'26.62"
$wsNew.Range("F2").Value = "'1.32"
$wsNew.Range("G2").Value = "'0.0180"
$wsNew.Range("H2").Value = 6

$wsNew.Range("A3").Value = 1
$wsNew.Range("B3").Value = "'011730"
$wsNew.Range("C3").Value = "工银聚享混合C"
$wsNew.Range("D3").Value = "'0.00"
$wsNew.Range("E3").Value = "'26.62"
$wsNew.Range("F3").Value = "'1.32"
$wsNew.Range("G3").Value = 0
$wsNew.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: push the existing two rows down
#    by one and insert the new "2022-Q3" summary row at the top.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Range("A3:D3").Copy($wsTotal.Range("A4"))
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3"))

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.02

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

# ---------------------------------------------------------------------
# 4. Restore "2022-Q1" (now the 4th tab) as the selected/active sheet,
#    matching the workbook's pre-existing active-tab state (adding the
#    new sheet would otherwise leave the freshly inserted sheet active).
# ---------------------------------------------------------------------
$wb.Worksheets.Item(4).Activate()
